# "Fruta / hortaliza, semanal"
# A new weekly price record was inserted at row 221 of the "Zanahoria"
# (carrot) price sheet for Terminal Hortofrutícola Agro Chillán, pushing
# every following record down by one row (221-270 -> 222-271).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 221, shifting rows 221:270 down to 222:271.
$ws.Rows("221:221").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(221, 1).Value = 7
$ws.Cells.Item(221, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(221, 3).Value = "Ñuble"
$ws.Cells.Item(221, 4).Value = 44637
$ws.Cells.Item(221, 5).Value = 16
$ws.Cells.Item(221, 6).Value = 100114013
$ws.Cells.Item(221, 7).Value = "Zanahoria"
$ws.Cells.Item(221, 8).Value = "Sin especificar"
$ws.Cells.Item(221, 9).Value = "Primera"
$ws.Cells.Item(221, 10).Value = 120
$ws.Cells.Item(221, 11).Value = 6000
$ws.Cells.Item(221, 12).Value = 6500
$ws.Cells.Item(221, 13).Value = 6250
$ws.Cells.Item(221, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(221, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(221, 16).Value = 312
$ws.Cells.Item(221, 17).Value = 20
$ws.Cells.Item(221, 18).Value = "Hortaliza"
